$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 20,10

$arr[0,0] = -19.4802340341602
$arr[0,1] = 2.41600707782063
$arr[0,2] = -19.4802340341602
$arr[0,3] = -19.4802340341602
$arr[0,4] = -19.4802340341602
$arr[0,5] = -19.4802340341602
$arr[0,6] = -19.4802340341602
$arr[0,7] = -19.4802340341602
$arr[0,8] = -19.4802340341602
$arr[0,9] = -19.4802340341602
$arr[1,0] = -19.4802340341602
$arr[1,1] = -19.4802340341602
$arr[1,2] = -19.4802340341602
$arr[1,3] = -19.4802340341602
$arr[1,4] = -19.4802340341602
$arr[1,5] = -19.4802340341602
$arr[1,6] = -19.4802340341602
$arr[1,7] = 2.317950187151309
$arr[1,8] = -19.4802340341602
$arr[1,9] = -19.4802340341602
$arr[2,0] = -19.4802340341602
$arr[2,1] = 2.162759397820151
$arr[2,2] = 2.820974383070498
$arr[2,3] = -19.4802340341602
$arr[2,4] = 2.545456063315606
$arr[2,5] = -19.4802340341602
$arr[2,6] = 2.077400275094933
$arr[2,7] = -19.4802340341602
$arr[2,8] = 4.321926220912938
$arr[2,9] = -19.4802340341602
$arr[3,0] = -19.4802340341602
$arr[3,1] = 0.992652111868738
$arr[3,2] = -19.4802340341602
$arr[3,3] = -19.4802340341602
$arr[3,4] = -19.4802340341602
$arr[3,5] = 1.84715468789251
$arr[3,6] = -19.4802340341602
$arr[3,7] = -19.4802340341602
$arr[3,8] = -19.4802340341602
$arr[3,9] = -19.4802340341602
$arr[4,0] = -19.4802340341602
$arr[4,1] = -19.4802340341602
$arr[4,2] = -19.4802340341602
$arr[4,3] = -19.4802340341602
$arr[4,4] = -19.4802340341602
$arr[4,5] = -19.4802340341602
$arr[4,6] = -19.4802340341602
$arr[4,7] = -19.4802340341602
$arr[4,8] = -19.4802340341602
$arr[4,9] = -19.4802340341602
$arr[5,0] = 2.95919537927759
$arr[5,1] = -19.4802340341602
$arr[5,2] = -19.4802340341602
$arr[5,3] = -19.4802340341602
$arr[5,4] = -19.4802340341602
$arr[5,5] = -19.4802340341602
$arr[5,6] = -19.4802340341602
$arr[5,7] = -19.4802340341602
$arr[5,8] = -19.4802340341602
$arr[5,9] = -19.4802340341602
$arr[6,0] = -19.4802340341602
$arr[6,1] = -19.4802340341602
$arr[6,2] = -19.4802340341602
$arr[6,3] = 2.844009427223499
$arr[6,4] = -19.4802340341602
$arr[6,5] = -19.4802340341602
$arr[6,6] = -19.4802340341602
$arr[6,7] = -19.4802340341602
$arr[6,8] = -19.4802340341602
$arr[6,9] = -19.4802340341602
$arr[7,0] = 3.61153529341555
$arr[7,1] = -19.4802340341602
$arr[7,2] = -19.4802340341602
$arr[7,3] = -19.4802340341602
$arr[7,4] = -19.4802340341602
$arr[7,5] = -19.4802340341602
$arr[7,6] = -19.4802340341602
$arr[7,7] = -19.4802340341602
$arr[7,8] = -19.4802340341602
$arr[7,9] = -19.4802340341602
$arr[8,0] = -19.4802340341602
$arr[8,1] = -19.4802340341602
$arr[8,2] = -19.4802340341602
$arr[8,3] = -19.4802340341602
$arr[8,4] = -19.4802340341602
$arr[8,5] = -19.4802340341602
$arr[8,6] = -19.4802340341602
$arr[8,7] = 1.380808180803917
$arr[8,8] = -19.4802340341602
$arr[8,9] = 2.439092327797603
$arr[9,0] = -19.4802340341602
$arr[9,1] = -19.4802340341602
$arr[9,2] = -19.4802340341602
$arr[9,3] = 2.027508836441826
$arr[9,4] = -19.4802340341602
$arr[9,5] = 2.743515821656592
$arr[9,6] = -19.4802340341602
$arr[9,7] = -19.4802340341602
$arr[9,8] = -19.4802340341602
$arr[9,9] = 1.28306668905348
$arr[10,0] = -19.4802340341602
$arr[10,1] = -19.4802340341602
$arr[10,2] = -19.4802340341602
$arr[10,3] = -19.4802340341602
$arr[10,4] = -19.4802340341602
$arr[10,5] = -19.4802340341602
$arr[10,6] = -19.4802340341602
$arr[10,7] = -19.4802340341602
$arr[10,8] = -19.4802340341602
$arr[10,9] = -19.4802340341602
$arr[11,0] = -19.4802340341602
$arr[11,1] = -19.4802340341602
$arr[11,2] = -19.4802340341602
$arr[11,3] = 1.740437731541907
$arr[11,4] = -19.4802340341602
$arr[11,5] = -19.4802340341602
$arr[11,6] = -19.4802340341602
$arr[11,7] = -19.4802340341602
$arr[11,8] = -19.4802340341602
$arr[11,9] = 1.630400287654835
$arr[12,0] = -19.4802340341602
$arr[12,1] = -19.4802340341602
$arr[12,2] = 1.716820040335221
$arr[12,3] = -19.4802340341602
$arr[12,4] = -19.4802340341602
$arr[12,5] = -19.4802340341602
$arr[12,6] = -19.4802340341602
$arr[12,7] = -19.4802340341602
$arr[12,8] = -19.4802340341602
$arr[12,9] = 2.224841864093364
$arr[13,0] = -19.4802340341602
$arr[13,1] = -19.4802340341602
$arr[13,2] = -0.1800848319395936
$arr[13,3] = -19.4802340341602
$arr[13,4] = -19.4802340341602
$arr[13,5] = -19.4802340341602
$arr[13,6] = -19.4802340341602
$arr[13,7] = -19.4802340341602
$arr[13,8] = -19.4802340341602
$arr[13,9] = -19.4802340341602
$arr[14,0] = -19.4802340341602
$arr[14,1] = -19.4802340341602
$arr[14,2] = -19.4802340341602
$arr[14,3] = -19.4802340341602
$arr[14,4] = -19.4802340341602
$arr[14,5] = -19.4802340341602
$arr[14,6] = -19.4802340341602
$arr[14,7] = -19.4802340341602
$arr[14,8] = -19.4802340341602
$arr[14,9] = -19.4802340341602
$arr[15,0] = -19.4802340341602
$arr[15,1] = 0.8378817450149565
$arr[15,2] = 0.08114528496325832
$arr[15,3] = -19.4802340341602
$arr[15,4] = -19.4802340341602
$arr[15,5] = -19.4802340341602
$arr[15,6] = 0.7362525191106877
$arr[15,7] = 1.293775196632285
$arr[15,8] = -19.4802340341602
$arr[15,9] = -19.4802340341602
$arr[16,0] = -19.4802340341602
$arr[16,1] = -19.4802340341602
$arr[16,2] = -19.4802340341602
$arr[16,3] = -19.4802340341602
$arr[16,4] = -19.4802340341602
$arr[16,5] = -19.4802340341602
$arr[16,6] = 0.4575262607533383
$arr[16,7] = 1.052004875937155
$arr[16,8] = -19.4802340341602
$arr[16,9] = -19.4802340341602
$arr[17,0] = -19.4802340341602
$arr[17,1] = -19.4802340341602
$arr[17,2] = 1.73382237978093
$arr[17,3] = -19.4802340341602
$arr[17,4] = -19.4802340341602
$arr[17,5] = -19.4802340341602
$arr[17,6] = 2.187922526626506
$arr[17,7] = 2.143181771134634
$arr[17,8] = -19.4802340341602
$arr[17,9] = -19.4802340341602
$arr[18,0] = -19.4802340341602
$arr[18,1] = 1.654373398740997
$arr[18,2] = 2.13089390330602
$arr[18,3] = -19.4802340341602
$arr[18,4] = 3.823958090109333
$arr[18,5] = -19.4802340341602
$arr[18,6] = 1.964785359741539
$arr[18,7] = 1.793753310466423
$arr[18,8] = -19.4802340341602
$arr[18,9] = 2.12856808108556
$arr[19,0] = -19.4802340341602
$arr[19,1] = 1.704973519701118
$arr[19,2] = -19.4802340341602
$arr[19,3] = 2.433355553923676
$arr[19,4] = -19.4802340341602
$arr[19,5] = 3.278739182645221
$arr[19,6] = 2.097764964856446
$arr[19,7] = -19.4802340341602
$arr[19,8] = -19.4802340341602
$arr[19,9] = -19.4802340341602

$ws.Range("B2:K21").Value = $arr
